# Edit: S4Loader fixed a primary key error on meta_pdt_field_value_enums table.
# Adds a new "PDT_Tank Data" worksheet (between "PDT_Pressure Transmitter" and
# "PDT_Flow Transmitter Data") populated with tank-specific PDT field
# specifications, and updates sheet selections/active-tab bookkeeping.

$wb = $excel.ActiveWorkbook

# --- Update selection on "PDT_Pressure Transmitter" (no longer the active tab) ---
$wsPressure = $wb.Worksheets.Item("PDT_Pressure Transmitter")
$wsPressure.Activate()
$wsPressure.Range("B2:F8").Select()

# --- Insert the new "PDT_Tank Data" sheet right after "PDT_Pressure Transmitter" ---
$newSheet = $wb.Worksheets.Add($null, $wsPressure)
$newSheet.Name = "PDT_Tank Data"
$ws = $newSheet

# Column widths (characters) to match the other PDT_* sheets.
$ws.Columns.Item(1).ColumnWidth = 41.5
$ws.Columns.Item(3).ColumnWidth = 15.67
$ws.Columns.Item(4).ColumnWidth = 25.5
$ws.Columns.Item(5).ColumnWidth = 25.5
$ws.Columns.Item(6).ColumnWidth = 25.5

# --- Populate the field-specification table ---
$ws.Range("A1").Value = 'pdt_field_name'
$ws.Range("B1").Value = 's4_type'
$ws.Range("C1").Value = 'evaluator_type'
$ws.Range("D1").Value = 'evaluator'
$ws.Range("E1").Value = 'builtin_arg1'
$ws.Range("F1").Value = 'sample'
$ws.Range("A2").Value = 'Asset Name'
$ws.Range("B2").Value = 'String'
$ws.Range("C2").Value = 'builtin'
$ws.Range("D2").Value = 'non-empty-string'
$ws.Range("F2").Value = 'DMP OAC 01 01 Damper-1'
$ws.Range("A3").Value = 'Asset Type'
$ws.Range("B3").Value = 'String'
$ws.Range("C3").Value = 'enum'
$ws.Range("A4").Value = 'Asset Status'
$ws.Range("B4").Value = 'String'
$ws.Range("C4").Value = 'builtin'
$ws.Range("D4").Value = 'asset-status'
$ws.Range("A5").Value = 'Uniclass Code'
$ws.Range("A6").Value = 'Uniclass Description'
$ws.Range("A7").Value = 'Manufacturer Data'
$ws.Range("A8").Value = 'Manufacturer'
$ws.Range("B8").Value = 'String'
$ws.Range("C8").Value = 'builtin'
$ws.Range("D8").Value = 'manufacturer'
$ws.Range("A9").Value = 'Manufacturer website'
$ws.Range("A10").Value = 'Product range'
$ws.Range("A11").Value = 'Product model number '
$ws.Range("A12").Value = 'CE approval '
$ws.Range("A13").Value = 'Approvals '
$ws.Range("A14").Value = 'Product literature'
$ws.Range("A15").Value = 'Features'
$ws.Range("A16").Value = 'Sustainability'
$ws.Range("A17").Value = 'Embodied carbon '
$ws.Range("A18").Value = 'Life cycle analysis'
$ws.Range("A19").Value = 'Location of manufacturer'
$ws.Range("A20").Value = 'Green guide for specification'
$ws.Range("A21").Value = 'Environmental product declaration'
$ws.Range("A22").Value = 'Responsible sourcing of materials'
$ws.Range("A23").Value = 'Energy technology list'
$ws.Range("A24").Value = 'Responsible extraction of materials'
$ws.Range("A25").Value = 'Material ingredient reporting'
$ws.Range("A26").Value = 'Operations & Maintenance'
$ws.Range("A27").Value = 'O & M manual'
$ws.Range("A28").Value = 'Daily'
$ws.Range("A29").Value = 'Weekly'
$ws.Range("A30").Value = 'Monthly'
$ws.Range("A31").Value = 'Quarterly'
$ws.Range("A32").Value = '6 Monthly'
$ws.Range("A33").Value = 'Annually'
$ws.Range("A34").Value = 'Bespoke timeframe'
$ws.Range("A35").Value = 'Maintenance required 0-300hrs'
$ws.Range("A36").Value = 'Maintenance required 301-600hrs'
$ws.Range("A37").Value = 'Maintenance required 601-1000hrs'
$ws.Range("A38").Value = 'Maintenance required 1001-2000hrs'
$ws.Range("A39").Value = 'Maintenance required 2001-4000hrs'
$ws.Range("A40").Value = 'Maintenance required 4001-8000hrs'
$ws.Range("A41").Value = 'Maintenance required 8001-12000hrs'
$ws.Range("A42").Value = 'Expected life'
$ws.Range("A43").Value = 'Warranty ID'
$ws.Range("A44").Value = 'Construction Data'
$ws.Range("A45").Value = 'Column Material'
$ws.Range("B45").Value = 'String'
$ws.Range("C45").Value = 'enum'
$ws.Range("A46").Value = 'Floor Material'
$ws.Range("B46").Value = 'String'
$ws.Range("C46").Value = 'enum'
$ws.Range("A47").Value = 'Reservoir Covering'
$ws.Range("B47").Value = 'String'
$ws.Range("C47").Value = 'enum'
$ws.Range("A48").Value = 'Roof Material'
$ws.Range("B48").Value = 'String'
$ws.Range("C48").Value = 'enum'
$ws.Range("A49").Value = 'Tank Construction'
$ws.Range("B49").Value = 'String'
$ws.Range("C49").Value = 'enum'
$ws.Range("A50").Value = 'Tank Covering'
$ws.Range("B50").Value = 'String'
$ws.Range("C50").Value = 'enum'
$ws.Range("A51").Value = 'Tank Level'
$ws.Range("B51").Value = 'String'
$ws.Range("C51").Value = 'enum'
$ws.Range("A52").Value = 'Wall Material'
$ws.Range("B52").Value = 'String'
$ws.Range("C52").Value = 'enum'
$ws.Range("A53").Value = 'Dimensions & Weight Data'
$ws.Range("A54").Value = 'Bottom Water Level (mAOD)'
$ws.Range("A55").Value = 'Top Water Level (mAOD)'
$ws.Range("A56").Value = 'Tank Shape'
$ws.Range("B56").Value = 'String'
$ws.Range("C56").Value = 'enum'
$ws.Range("A57").Value = 'Capacity (m3)'
$ws.Range("B57").Value = 'Decimal'
$ws.Range("C57").Value = 'builtin'
$ws.Range("D57").Value = 'decimal-with-units'
$ws.Range("E57").Value = 'm3'
$ws.Range("F57").Value = '50.4 m3'
$ws.Range("A58").Value = 'Diameter (mm)'
$ws.Range("B58").Value = 'Int'
$ws.Range("C58").Value = 'builtin'
$ws.Range("D58").Value = 'integer-with-units'
$ws.Range("E58").Value = 'mm'
$ws.Range("F58").Value = '800 mm'
$ws.Range("A59").Value = 'Side Depth (mm)'
$ws.Range("B59").Value = 'Int'
$ws.Range("C59").Value = 'builtin'
$ws.Range("D59").Value = 'integer-with-units'
$ws.Range("E59").Value = 'mm'
$ws.Range("F59").Value = '800 mm'
$ws.Range("A60").Value = 'Side Depth Max (mm)'
$ws.Range("B60").Value = 'Int'
$ws.Range("C60").Value = 'builtin'
$ws.Range("D60").Value = 'integer-with-units'
$ws.Range("E60").Value = 'mm'
$ws.Range("F60").Value = '800 mm'
$ws.Range("A61").Value = 'Side Depth Min (mm)'
$ws.Range("B61").Value = 'Int'
$ws.Range("C61").Value = 'builtin'
$ws.Range("D61").Value = 'integer-with-units'
$ws.Range("E61").Value = 'mm'
$ws.Range("F61").Value = '4000 mm'
$ws.Range("A62").Value = 'Top Surface Area (m2)'
$ws.Range("B62").Value = 'Decimal'
$ws.Range("C62").Value = 'builtin'
$ws.Range("D62").Value = 'decimal-with-units'
$ws.Range("E62").Value = 'm2'
$ws.Range("F62").Value = '35.4 m2'
$ws.Range("A63").Value = 'Centre Depth (mm)'
$ws.Range("B63").Value = 'Int'
$ws.Range("C63").Value = 'builtin'
$ws.Range("D63").Value = 'integer-with-units'
$ws.Range("E63").Value = 'mm'
$ws.Range("F63").Value = '4000 mm'
$ws.Range("A64").Value = 'Length (mm)'
$ws.Range("B64").Value = 'Int'
$ws.Range("C64").Value = 'builtin'
$ws.Range("D64").Value = 'integer-with-units'
$ws.Range("E64").Value = 'mm'
$ws.Range("F64").Value = '4000 mm'
$ws.Range("A65").Value = 'Major Axis (mm)'
$ws.Range("A66").Value = 'Minor Axis (mm)'
$ws.Range("A67").Value = 'Asset Data '
$ws.Range("A68").Value = 'Manufacturer''s serial number'
$ws.Range("B68").Value = 'String'
$ws.Range("C68").Value = 'builtin'
$ws.Range("D68").Value = 'non-empty-string'
$ws.Range("A69").Value = 'Date of installation'
$ws.Range("B69").Value = 'Date'
$ws.Range("C69").Value = 'builtin'
$ws.Range("D69").Value = 'local-date'
$ws.Range("A70").Value = 'Tag reference'
$ws.Range("B70").Value = 'String'
$ws.Range("C70").Value = 'builtin'
$ws.Range("D70").Value = 'any'
$ws.Range("A71").Value = 'Legacy GUID (Globally Unique Identifier) reference'
$ws.Range("B71").Value = 'String'
$ws.Range("C71").Value = 'builtin'
$ws.Range("D71").Value = 'any'
$ws.Range("A72").Value = 'Location in Asset Hierarchy'
$ws.Range("B72").Value = 'String'
$ws.Range("C72").Value = 'builtin'
$ws.Range("D72").Value = 'floc'
$ws.Range("F72").Value = 'ABB01-SSS-SFS-HSE-SYS02'
$ws.Range("A73").Value = 'Specific Model'
$ws.Range("B73").Value = 'String'
$ws.Range("C73").Value = 'builtin'
$ws.Range("D73").Value = 'non-empty-string'
$ws.Range("A74").Value = 'Location on Site'
$ws.Range("B74").Value = 'String'
$ws.Range("C74").Value = 'builtin'
$ws.Range("D74").Value = 'any'
$ws.Range("A75").Value = 'Manufacturers Asset Life (yr)'
$ws.Range("A76").Value = 'Memo Line'
$ws.Range("B76").Value = 'String'
$ws.Range("C76").Value = 'builtin'
$ws.Range("D76").Value = 'any'
$ws.Range("A77").Value = 'Content Type'
$ws.Range("B77").Value = 'String'
$ws.Range("C77").Value = 'enum'

# --- Update the "PDT_Flow Transmitter Data" sheet's scroll position / selection ---
$wsFlow = $wb.Worksheets.Item("PDT_Flow Transmitter Data")
$wsFlow.Activate()
$wsFlow.Range("B61:F65").Select()
$excel.ActiveWindow.ScrollRow = 41

# --- Make the newly added sheet the active tab/selection, as in the final workbook ---
$newSheet.Activate()
$newSheet.Range("B1").Select()
